$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '96.878.81'
$ws.Range('E2').Value = '  -1.08%  '
$ws.Range('D3').Value = '3.340.70'
$ws.Range('E3').Value = '  -2.47%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '250.44'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.87%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '656.22'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.02%  '
$ws.Range('E7').Value = '  -4.99%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.424'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.17%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.999'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.03%  '
$ws.Range('E10').Value = '  -5.31%  '
$ws.Range('D11').Value = '3.337.43'
$ws.Range('E11').Value = '  -2.44%  '
$ws.Range('E12').Value = '  -2.65%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '40.78'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.78%  '
$ws.Range('D14').Value = '96.588.54'
$ws.Range('E14').Value = '  -1.07%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.10'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.72%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000253'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.51%  '
$ws.Range('D17').Value = '3.967.16'
$ws.Range('E17').Value = '  -2.57%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '8.73'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.49%  '
$ws.Range('D19').Value = '3.336.40'
$ws.Range('E19').Value = '  -2.36%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.564'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +12.98%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.52'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.11%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.70'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.35%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '509.23'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.70%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.34'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.63%  '
$ws.Range('E25').Value = '  -3.74%  '
$ws.Range('E26').Value = '  +6.34%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '96.66'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.49%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '12.12'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -5.14%  '
$ws.Range('E29').Value = '  -3.36%  '
$ws.Range('E30').Value = '  +0.41%  '
$ws.Range('E31').Value = '  -0.15%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.189'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.54'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +10.87%  '
$ws.Range('E34').Value = '  -0.05%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.554'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.10%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '28.42'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.73%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.50'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +4.56%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '7.82'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.07%  '
$ws.Range('E39').Value = '  +0.02%  '
$ws.Range('E40').Value = '  -2.23%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '507.53'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.79%  '
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0437'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.31%  '
$ws.Range('B43').Value = 'WhiteBITCoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '24.36'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.58%  '
$ws.Range('B44').Value = 'ARBITRUM'
$ws.Range('C44').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.838'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.04%  '
$ws.Range('B45').Value = 'MantraDAO'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.68'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.12%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.60'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.92%  '
$ws.Range('E47').Value = '  +5.05%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.52'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.01%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '54.86'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +6.89%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.11'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -6.23%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '162.18'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.39%  '
